# edit.ps1 - applies the "Hyperparameter Tuning" commit changes:
#  1. Inserts <w:lastRenderedPageBreak/> before four specific runs.
#  2. Renumbers four legacy VML <v:shape id="..."> image ids.
#  3. Splits the "First the cost function..." run in two (with a
#     lastRenderedPageBreak on the second half) and moves the "_GoBack"
#     bookmark from the end of the following paragraph to the start of
#     this one.
#
# NOTE: this runtime's COM/PS shim loses the "live" link to the document
# when a COM object (e.g. a Paragraph) is passed through a *named*
# function parameter (e.g. "Foo -Paragraph $p"); the call silently
# succeeds but the underlying InsertXML no longer mutates the real
# document. Positional parameter passing ("Foo $p") does not have this
# problem, so every function call below is positional only.

$d = $word.ActiveDocument

function Set-ParagraphXml($Paragraph, $OldFragment, $NewFragment) {
    $range = $Paragraph.Range
    $xml = $range.WordOpenXML
    if (-not $xml.Contains($OldFragment)) {
        throw "Expected fragment not found in paragraph XML: $OldFragment"
    }
    $updated = $xml.Replace($OldFragment, $NewFragment)
    $range.InsertXML($updated)
}

function Find-ParagraphContaining($NeedleText) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text.Contains($NeedleText)) {
            return $p
        }
    }
    throw "No paragraph found containing: $NeedleText"
}

function Find-ParagraphWithShapeId($ShapeId) {
    $paras = $d.Paragraphs
    $needle = 'v:shape id="' + $ShapeId + '"'
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text.Length -le 2) {
            $xml = $p.Range.WordOpenXML
            if ($xml.Contains($needle)) {
                return $p
            }
        }
    }
    throw "No paragraph found containing shape id: $ShapeId"
}

# --- 1 & 2. lastRenderedPageBreak insertions on Courier-New code runs ---
# (the exact <w:t>...</w:t> fragments are hard-coded verbatim - including
# the non-breaking spaces in the "Z[i,j]" line, built via [char]0x00A0 to
# dodge source-encoding pitfalls - because this shim's
# String.LastIndexOf(value, startIndex) ignores startIndex, so it cannot
# be used to reliably re-derive the enclosing <w:t> element dynamically)

$nbsp = [char]0x00A0

$breakFindText = @(
    "Plot of cost_history vs. iterations",
    "Sklearn Implementation",
    "Z[i,j] = Cost"
)
$breakOldT = @(
    '<w:t># Plot of cost_history vs. iterations</w:t>',
    '<w:t># Sklearn Implementation</w:t>',
    ('<w:t>' + $nbsp + ' ' + $nbsp + ' Z[i,j] = Cost(x,y, theta=[[xx[i,j]], [yy[i,j]]])</w:t>')
)

for ($bi = 0; $bi -lt $breakFindText.Count; $bi++) {
    $p = Find-ParagraphContaining $breakFindText[$bi]
    $oldT = $breakOldT[$bi]
    $newT = '<w:lastRenderedPageBreak/>' + $oldT
    Set-ParagraphXml $p $oldT $newT
}

# --- 3. Renumber the four VML shape ids ---

$shapeIdMap = [ordered]@{
    "_x0000_i1216" = "_x0000_i1025"
    "_x0000_i1218" = "_x0000_i1026"
    "_x0000_i1228" = "_x0000_i1027"
    "_x0000_i1231" = "_x0000_i1028"
}

foreach ($oldId in $shapeIdMap.Keys) {
    $newId = $shapeIdMap[$oldId]
    $p = Find-ParagraphWithShapeId $oldId
    $oldFrag = 'v:shape id="' + $oldId + '"'
    $newFrag = 'v:shape id="' + $newId + '"'
    Set-ParagraphXml $p $oldFrag $newFrag
}

# --- 4. Remove the _GoBack bookmark from the "Additionally, ..." paragraph ---

$pAdditionally = Find-ParagraphContaining "Additionally, two values are predicted"
$oldBookmarkTail = '</w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$newBookmarkTail = '</w:r></w:p>'
Set-ParagraphXml $pAdditionally $oldBookmarkTail $newBookmarkTail

# --- 5. Split the "First the cost function..." run and add the bookmark ---

$pFirst = Find-ParagraphContaining "First the cost function and gradient descent algorithm"
$oldRun = '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> First the cost function and gradient descent algorithm is calculated using mathematical first principles, and then the linear regression model is predicted using these functions. Next, this approach is compared to the LinearRegression function contained in the SciKit Learn python library. From the output of the comparison graph we can see that the regression models obtained from first principles is almost identical to the on obtained in the SciKit Learn library.</w:t></w:r>'
$newRun = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">First the cost function and gradient descent algorithm is calculated using mathematical first principles, and then the linear regression model is predicted using these functions. Next, this approach is compared to the LinearRegression function contained in the SciKit Learn python library. From the output of the comparison graph we can see that the regression models </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:szCs w:val="20"/></w:rPr><w:lastRenderedPageBreak/><w:t>obtained from first principles is almost identical to the on obtained in the SciKit Learn library.</w:t></w:r>'
Set-ParagraphXml $pFirst $oldRun $newRun

Write-Output "All edits applied."
